$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.042.66'
$ws.Range("E2").Value = '  +0.89%  '

$ws.Range("D3").Value = '2.650.14'
$ws.Range("E3").Value = '  +1.38%  '

$ws.Range("D5").Value = '''533.11'
$ws.Range("E5").Value = '  +4.32%  '

$ws.Range("D6").Value = '''155.99'
$ws.Range("E6").Value = '  +1.05%  '

$ws.Range("D7").Value = '''0.998'

$ws.Range("D8").Value = '''0.591'
$ws.Range("E8").Value = '  +0.71%  '

$ws.Range("D9").Value = '''6.61'
$ws.Range("E9").Value = '  -1.24%  '

$ws.Range("E10").Value = '  +5.04%  '

$ws.Range("E11").Value = '  +1.55%  '

$ws.Range("E12").Value = '  -0.15%  '

$ws.Range("D13").Value = '3.113.53'
$ws.Range("E13").Value = '  +1.39%  '

$ws.Range("D14").Value = '61.059.23'
$ws.Range("E14").Value = '  +1.00%  '

$ws.Range("D15").Value = '''22.09'

$ws.Range("E16").Value = '  +2.39%  '

$ws.Range("D17").Value = '2.648.10'
$ws.Range("E17").Value = '  +0.86%  '

$ws.Range("E18").Value = '  +0.13%  '

$ws.Range("D19").Value = '''354.99'
$ws.Range("E19").Value = '  +1.18%  '

$ws.Range("D20").Value = '''10.68'
$ws.Range("E20").Value = '  +0.51%  '

$ws.Range("D21").Value = '''6.25'
$ws.Range("E21").Value = '  +1.42%  '

$ws.Range("D22").Value = '''1.00'
$ws.Range("E22").Value = '  +0.51%  '

$ws.Range("D23").Value = '''61.68'
$ws.Range("E23").Value = '  +1.77%  '

$ws.Range("E24").Value = '  +2.08%  '

$ws.Range("E25").Value = '  +1.34%  '

$ws.Range("E26").Value = '  +0.37%  '

$ws.Range("D27").Value = '0.0₃0859'
$ws.Range("E27").Value = '  +1.81%  '

$ws.Range("E28").Value = '  +0.20%  '

$ws.Range("E29").Value = '  -0.02%  '

$ws.Range("E30").Value = '  +7.02%  '

$ws.Range("E31").Value = '  +4.12%  '

$ws.Range("E32").Value = '  +0.54%  '

$ws.Range("D33").Value = '''150.02'
$ws.Range("E33").Value = '  -0.30%  '

$ws.Range("E34").Value = '  +3.39%  '

$ws.Range("E35").Value = '  +1.04%  '

$ws.Range("D36").Value = '''0.922'
$ws.Range("E36").Value = '  +8.79%  '

$ws.Range("D37").Value = '''0.896'
$ws.Range("E37").Value = '  +1.75%  '

$ws.Range("D38").Value = '''309.24'
$ws.Range("E38").Value = '  +5.16%  '

$ws.Range("E39").Value = '  +0.75%  '

$ws.Range("D40").Value = '''3.82'
$ws.Range("E40").Value = '  +1.30%  '

$ws.Range("D41").Value = '''0.647'
$ws.Range("E41").Value = '  +3.36%  '

$ws.Range("E42").Value = '  +1.73%  '

$ws.Range("E43").Value = '  +1.41%  '

$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = '''0.998'
$ws.Range("E44").Value = '  +0.11%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '''19.97'
$ws.Range("E45").Value = '  +0.24%  '

$ws.Range("D46").Value = '''5.00'
$ws.Range("E46").Value = '  +2.16%  '

$ws.Range("D47").Value = '''0.0239'
$ws.Range("E47").Value = '  +2.25%  '

$ws.Range("D48").Value = '''19.22'
$ws.Range("E48").Value = '  +7.82%  '

$ws.Range("D49").Value = '''10.36'
$ws.Range("E49").Value = '  +0.36%  '

$ws.Range("D50").Value = '1.989.27'
$ws.Range("E50").Value = '  -0.69%  '

$ws.Range("D51").Value = '''1.84'
$ws.Range("E51").Value = '  +2.19%  '
